# Update the Book Chapter IDs for the two rows from the old "FT_" prefixed
# placeholder identifiers to the final "BK_" (Book) identifiers, and restore
# the sheet selection/scroll position that was active when the edit was made.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "BK_001"
$ws.Range("A3").Value = "BK_002"

# Move the selection to C3 and scroll the view back to the top-left (A1),
# matching the saved sheetView state.
$null = $ws.Range("C3").Select()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
